$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values like "29.552.67"
# or "326.39" are not auto-converted to numbers by Excel type inference.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) columns ---
$ws.Range("D2").Value = "29.552.67"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.926.39"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "326.39"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "0.4818"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.4056"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.08204"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "1.010"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "23.81"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "1.936.14"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "6.102"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "7.306"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "91.55"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "0.06889"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "17.67"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "29.562.56"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "5.676"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").Value = "2.176"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "2.172.19"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "155.87"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "6.410"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "2.095"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "120.64"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "1.014"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "0.09584"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "5.594"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").Value = "3.561"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "1.383"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "0.06357"
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "1.191"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "0.5954"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "10.70"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D43").Value = "0.1845"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "2.476"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "1.281"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "12.47"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "0.07478"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "0.5553"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "1.977"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "118.51"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").Value = "2.435"
$ws.Range("E51").Value = "  +1.19%  "

# Row 41 and 42 swap positions: "Frax" now appears before "FraxShare",
# each carrying its own updated link, price and volume figures.
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "1.011"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "7.895"
$ws.Range("E42").Value = "  -1.16%  "

# Restore default style on column D so only the number format applied
# above (which forced text storage) is cleared, matching the original
# unstyled appearance of these cells.
$ws.Range("D2:D51").Style = "Normal"
